$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row values ---
$ws.Cells.Item(1, 1).Value = "Common Word"
$ws.Cells.Item(1, 2).Value = "Total Frequency"
$ws.Cells.Item(1, 3).Value = "Websites"

# --- Data rows ---
$data = @(
    @("Cookies", 140, "https://www.100-dakar.com (14), https://www.benjaminwahl.at (18), https://www.dasholzhaus.at (11), https://www.diequote.at (24), https://www.drehorgelkabarett.at (14), https://www.frautomani.at (4), https://www.ingridschiller.at (14), https://www.ottosaxinger.at (3), https://www.peligro.at (14), https://www.reinhardreisenzahn.com (4), https://www.schuledesungehorsams.at (2), https://www.skodone.at (18)"),
    @("Linz", 55, "https://www.freie-medien.at (10), https://www.freizeitundkommunikation.at (3), https://www.linzfmr.at (18), https://www.pflueckt.at (3), https://www.steingeschichten.at (21)"),
    @("Film", 47, "https://www.apileofghosts.com (16), https://www.boxafilm.com (8), https://www.corpushomini.info (3), https://www.doublehappiness.at (12), https://www.retrogoldmine.com (8)"),
    @("Art", 45, "https://www.eipcp.net (13), https://www.faces-l.net (16), https://www.kairus.org (6), https://www.negentropy-sport.net (2), https://www.radical-openness.org (8)"),
    @("Page", 35, "https://www.hungaromedia.at (8), https://www.kuenstlerinnen.at (8), https://www.luckeneder-art.at (8), https://www.platform-socialism.org (3), https://www.regional-express.org (8)"),
    @("March", 49, "https://www.das-kollektiv.at (8), https://www.feminismus-krawall.at (16), https://www.fiftitu.at (19), https://www.unkraut-comics.at (6)"),
    @("School", 294, "https://www.alteschule-gutau.at (4), https://www.derschueler.at (5), https://www.die-schule.at (285)"),
    @("Culture", 36, "https://www.frauenkultur.at (17), https://www.igkultur.at (14), https://www.interregnum.live (5)"),
    @("Radio", 175, "https://www.radio-fri.at (6), https://www.schulradiotag.at (169)"),
    @("More", 64, "https://www.diebresche.org (23), https://www.programmkinowels.at (41)"),
    @("Uhr", 32, "https://www.fro.at (20), https://www.rudolfhabringer.at (12)"),
    @("Magdalena", 12, "https://www.magdalenareiter.at (2), https://www.themagdalenaproject.org (10)"),
    @("Casino", 10, "https://www.frf.at (6), https://www.photosalonhelga.com (4)"),
    @("Andreas", 8, "https://www.andreaskurz.net (2), https://www.andreaszingerle.com (6)"),
    @("Anna", 8, "https://www.anna-kraher.de (5), https://www.lllk.at (3)")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}

# --- Header formatting: bold, centered, top-aligned, thin box border ---
# Build the combined style on an out-of-the-way scratch cell first, then
# paste just the formats onto the header range in a single operation so the
# workbook only gains one new font / one new border / one new cell style
# (mirrors how the sheet was produced upstream) instead of one per property.
$scratch = $ws.Range("Z100")
$scratch.Font.Bold = $true
$scratch.HorizontalAlignment = -4108
$scratch.VerticalAlignment = -4160
$scratch.Borders.LineStyle = 1
$scratch.Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)
$scratch.Clear()
$ws.Range("A1").Select()

